# Updates cryptos list data (coin prices / 1h volume %) to the latest
# scrape, matching the GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.703.00'
$ws.Range('E2').Value = '  -0.51%  '

$ws.Range('D3').Value = '1.630.03'

$ws.Range('E4').Value = '  +0.44%  '

$ws.Range('D5').Value = "'214.33"
$ws.Range('E5').Value = '  -0.61%  '

$ws.Range('E6').Value = '  -1.05%  '

$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.32%  '

$ws.Range('E8').Value = '  -0.78%  '

$ws.Range('D9').Value = "'0.0636"
$ws.Range('E9').Value = '  -1.21%  '

$ws.Range('D10').Value = "'19.44"
$ws.Range('E10').Value = '  -5.54%  '

$ws.Range('D11').Value = "'0.0784"
$ws.Range('E11').Value = '  +0.14%  '

$ws.Range('D12').Value = '1.628.28'
$ws.Range('E12').Value = '  -0.55%  '

$ws.Range('E13').Value = '  -1.22%  '

$ws.Range('D14').Value = '1.854.96'
$ws.Range('E14').Value = '  -0.60%  '

$ws.Range('D15').Value = "'0.550"
$ws.Range('E15').Value = '  -2.12%  '

$ws.Range('D16').Value = '0.0₃0767'
$ws.Range('E16').Value = '  -0.59%  '

$ws.Range('D17').Value = "'63.12"
$ws.Range('E17').Value = '  -0.17%  '

$ws.Range('D18').Value = '25.736.78'
$ws.Range('E18').Value = '  -0.47%  '

$ws.Range('E19').Value = '  +0.30%  '

$ws.Range('E20').Value = '  +1.11%  '

$ws.Range('D21').Value = "'193.69"
$ws.Range('E21').Value = '  -0.06%  '

$ws.Range('D22').Value = "'9.93"
$ws.Range('E22').Value = '  -0.27%  '

$ws.Range('D23').Value = "'6.20"
$ws.Range('E23').Value = '  +1.00%  '

$ws.Range('E24').Value = '  +0.41%  '

$ws.Range('D25').Value = "'1.78"
$ws.Range('E25').Value = '  -0.77%  '

$ws.Range('D26').Value = "'140.45"
$ws.Range('E26').Value = '  -0.03%  '

$ws.Range('D27').Value = "'0.119"
$ws.Range('E27').Value = '  -3.88%  '

$ws.Range('D28').Value = "'6.80"
$ws.Range('E28').Value = '  -0.64%  '

$ws.Range('E29').Value = '  -0.47%  '

$ws.Range('E30').Value = '  -0.74%  '

$ws.Range('D31').Value = "'0.0483"
$ws.Range('E31').Value = '  -3.01%  '

$ws.Range('D32').Value = "'3.33"
$ws.Range('E32').Value = '  +0.76%  '

$ws.Range('D33').Value = "'3.24"
$ws.Range('E33').Value = '  -0.06%  '

$ws.Range('D34').Value = "'1.58"
$ws.Range('E34').Value = '  -0.08%  '

$ws.Range('E35').Value = '  +0.53%  '

$ws.Range('D36').Value = "'0.893"
$ws.Range('E36').Value = '  -1.64%  '

$ws.Range('E37').Value = '  -0.32%  '

$ws.Range('E38').Value = '  -1.89%  '

$ws.Range('D39').Value = '1.104.04'
$ws.Range('E39').Value = '  -2.28%  '

$ws.Range('D40').Value = "'0.0155"
$ws.Range('E40').Value = '  -0.94%  '

$ws.Range('E41').Value = '  +0.77%  '

$ws.Range('D42').Value = "'5.57"
$ws.Range('E42').Value = '  +0.63%  '

$ws.Range('E43').Value = '  +1.37%  '

$ws.Range('D44').Value = "'0.798"
$ws.Range('E44').Value = '  -0.56%  '

$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.763.56'
$ws.Range('E45').Value = '  -0.81%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0110'
$ws.Range('E46').Value = '  -0.50%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'54.97"
$ws.Range('E47').Value = '  -1.42%  '

$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'0.419"
$ws.Range('E48').Value = '  -1.90%  '

$ws.Range('E49').Value = '  -0.22%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'7.66"
$ws.Range('E50').Value = '  -1.59%  '

$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = "'2.32"
$ws.Range('E51').Value = '  +3.82%  '

$wb.Save()
